# This script applies the inventory adjustment changes described in the
# commit diff for CryCompanywiseStockReport_1.xlsx.
#
# For a set of stock rows, the quantity (column F) and stock value
# (column G = Rate x Qty) are corrected, and the "Sub Total:" rows
# (column B) for the affected company groups are updated to match the
# corrected sum of stock values within each group. The final aggregate
# "Sub Total:" (B619) and "Grand Total:" (B620) rows are also updated
# to reflect the new grand total.
#
# Two pairs of rows (227/228 and 572/573) had their code/rate/qty/value
# data effectively swapped between the two rows in the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 74
$ws.Range("G6").Value = 2211.12
$ws.Range("B10").Value = 28226.01
$ws.Range("F16").Value = 60
$ws.Range("G16").Value = 3073.8
$ws.Range("B19").Value = 3175.49
$ws.Range("F71").Value = 321
$ws.Range("G71").Value = 20447.7
$ws.Range("F77").Value = 250
$ws.Range("G77").Value = 11685
$ws.Range("B90").Value = 177651.8
$ws.Range("F205").Value = 20
$ws.Range("G205").Value = 7542.8
$ws.Range("B216").Value = 39390.91
$ws.Range("F222").Value = 12
$ws.Range("G222").Value = 1739.16
$ws.Range("F223").Value = 12
$ws.Range("G223").Value = 1589.76
$ws.Range("B227").Value = 63520
$ws.Range("E227").Value = 153.4
$ws.Range("F227").Value = 66
$ws.Range("G227").Value = 9522.48
$ws.Range("B228").Value = 55373
$ws.Range("E228").Value = 163.62
$ws.Range("F228").Value = -94
$ws.Range("G228").Value = -13562.32
$ws.Range("F250").Value = 8
$ws.Range("G250").Value = 3955.04
$ws.Range("F255").Value = 553
$ws.Range("G255").Value = 94745.49000000001
$ws.Range("B260").Value = 184932.01
$ws.Range("F278").Value = 10
$ws.Range("G278").Value = 1373.2
$ws.Range("F303").Value = 29
$ws.Range("G303").Value = 6115.81
$ws.Range("B304").Value = 172453.19
$ws.Range("F328").Value = 41
$ws.Range("G328").Value = 1525.61
$ws.Range("B330").Value = 27360.5
$ws.Range("F354").Value = 16
$ws.Range("G354").Value = 1097.44
$ws.Range("B358").Value = 35417.94
$ws.Range("F422").Value = 12
$ws.Range("G422").Value = 2593.32
$ws.Range("B424").Value = 3025.99
$ws.Range("F434").Value = 6
$ws.Range("G434").Value = 195.84
$ws.Range("B435").Value = 341.32
$ws.Range("F440").Value = 1
$ws.Range("G440").Value = 321.89
$ws.Range("B445").Value = 7053.27
$ws.Range("F450").Value = 9
$ws.Range("G450").Value = 1248.66
$ws.Range("B460").Value = 13268.33
$ws.Range("F462").Value = 95
$ws.Range("G462").Value = 3118.85
$ws.Range("B475").Value = 45321.05
$ws.Range("F477").Value = 8
$ws.Range("G477").Value = 362.72
$ws.Range("B478").Value = 362.72
$ws.Range("F509").Value = 213
$ws.Range("G509").Value = 17120.94
$ws.Range("B510").Value = 22837.64
$ws.Range("F552").Value = 15
$ws.Range("G552").Value = 1526.85
$ws.Range("B560").Value = 4405.61
$ws.Range("B572").Value = 65362
$ws.Range("F572").Value = 20
$ws.Range("G572").Value = 817.4
$ws.Range("B573").Value = 65079
$ws.Range("F573").Value = 6
$ws.Range("G573").Value = 245.22
$ws.Range("F578").Value = 76
$ws.Range("G578").Value = 3791.64
$ws.Range("F580").Value = 54
$ws.Range("G580").Value = 3077.46
$ws.Range("B583").Value = 15812.42
$ws.Range("F599").Value = 1589
$ws.Range("G599").Value = 259181.79
$ws.Range("F601").Value = 401
$ws.Range("G601").Value = 113430.87
$ws.Range("B606").Value = 421050.56
$ws.Range("F612").Value = 30
$ws.Range("G612").Value = 1229.7
$ws.Range("B618").Value = 43448.72
$ws.Range("B619").Value = 1720744.51
$ws.Range("B620").Value = 1720744.51
